# Generate Report for Handoff
# Adds a new row (for the newly handed-off file "e6b13d2b-8398-...md")
# to the Overview, zh-cn and de-de tables.

$wb = $excel.ActiveWorkbook

$newMd        = "e6b13d2b-8398-4d27-a770-af7ae6d1b59dooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newPathMd    = "e2e\" + $newMd
$newZhCnXlf   = "e6b13d2b-8398-4d27-a770-af7ae6d1b59doooooooooooooooooooooooooooooooooooooooo.9bea847d045fc0dfd3249275264b7879636a990a.zh-cn.xlf"
$newDeDeXlf   = "e6b13d2b-8398-4d27-a770-af7ae6d1b59doooooooooooooooooooooooooooooooooooooooo.9bea847d045fc0dfd3249275264b7879636a990a.de-de.xlf"
$newUrl       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a7a2cf80d5c02f24c0e60f85b6d1b022394cd9f5/e2e/$newMd"

$hoDateZhCn   = "2016-09-02 02:33:49"
$hoDateDeDe   = "2016-09-02 02:33:53"
$statusReady  = "Ready for handoff"
$hbDateEmpty  = "0001-01-01 00:00:00"

function Style-LikeHyperlink($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Underline = $true
    $range.Font.ThemeColor = 0
    $range.Font.Color = 15570276
}

function Style-LikeDate($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ===================== Overview sheet =====================
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A3").Value = $newMd
$wsOv.Range("B3").Value = $newPathMd
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $newUrl, "", "", $newPathMd) | Out-Null
Style-LikeHyperlink $wsOv.Range("B3")
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("D3").Value = ""
$wsOv.Range("E3").Value = $statusReady
$wsOv.Range("F3").Value = $statusReady
$wsOv.Range("G3").Value = $hoDateDeDe
Style-LikeDate $wsOv.Range("G3")

$wsOv.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOv.Columns.Item(6).ColumnWidth = 17.2159881591797

# ===================== zh-cn sheet =====================
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $newMd
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newUrl, "", "", $newMd) | Out-Null
Style-LikeHyperlink $wsZh.Range("A3")
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $statusReady
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $newZhCnXlf
$wsZh.Range("H3").Value = $hoDateZhCn
Style-LikeDate $wsZh.Range("H3")
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = $hbDateEmpty
Style-LikeDate $wsZh.Range("K3")
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Columns.Item(3).ColumnWidth = 17.2159881591797

# ===================== de-de sheet =====================
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $newMd
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newUrl, "", "", $newMd) | Out-Null
Style-LikeHyperlink $wsDe.Range("A3")
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $statusReady
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $newDeDeXlf
$wsDe.Range("H3").Value = $hoDateDeDe
Style-LikeDate $wsDe.Range("H3")
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = $hbDateEmpty
Style-LikeDate $wsDe.Range("K3")
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Columns.Item(3).ColumnWidth = 17.2159881591797

Write-Host "Report generated for handoff"
